$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.807.03"
$ws.Range("E2").Value = "  -2.84%  "

$ws.Range("D3").Value = "1.615.18"
$ws.Range("E3").Value = "  -3.32%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").Value = "'306.56"
$ws.Range("E6").Value = "  -1.89%  "

$ws.Range("D7").Value = "'0.3896"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'0.3793"
$ws.Range("E8").Value = "  -3.48%  "

$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("D11").Value = "'48.52"
$ws.Range("E11").Value = "  -6.00%  "

$ws.Range("D12").Value = "'0.08424"
$ws.Range("E12").Value = "  -2.32%  "

$ws.Range("D13").Value = "'23.71"
$ws.Range("E13").Value = "  -6.37%  "

$ws.Range("D14").Value = "'6.990"

$ws.Range("D15").Value = "'0.00001268"
$ws.Range("E15").Value = "  -3.82%  "

$ws.Range("D16").Value = "'7.406"
$ws.Range("E16").Value = "  -3.81%  "

$ws.Range("D17").Value = "1.613.43"
$ws.Range("E17").Value = "  -3.61%  "

$ws.Range("D18").Value = "'93.18"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").Value = "'0.06907"
$ws.Range("E19").Value = "  -1.83%  "

$ws.Range("D20").Value = "'19.89"
$ws.Range("E20").Value = "  -4.65%  "

$ws.Range("D21").Value = "'6.772"
$ws.Range("E21").Value = "  -3.99%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "'13.33"
$ws.Range("E23").Value = "  -4.27%  "

$ws.Range("D24").Value = "23.787.39"
$ws.Range("E24").Value = "  -2.88%  "

$ws.Range("D25").Value = "'2.405"
$ws.Range("E25").Value = "  +1.54%  "

$ws.Range("D26").Value = "'2.804"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("D27").Value = "'22.02"
$ws.Range("E27").Value = "  -4.81%  "

$ws.Range("D28").Value = "'156.90"

$ws.Range("D29").Value = "'138.86"
$ws.Range("E29").Value = "  -5.78%  "

$ws.Range("D30").Value = "'5.240"
$ws.Range("E30").Value = "  -10.96%  "

$ws.Range("D31").Value = "'7.724"
$ws.Range("E31").Value = "  -6.38%  "

$ws.Range("D32").Value = "'2.480"
$ws.Range("E32").Value = "  -1.12%  "

$ws.Range("D33").Value = "1.791.04"
$ws.Range("E33").Value = "  -3.55%  "

$ws.Range("D34").Value = "'0.08049"
$ws.Range("E34").Value = "  -3.58%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02850"
$ws.Range("E35").Value = "  -5.60%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.9479"
$ws.Range("E36").Value = "  -3.54%  "

$ws.Range("D37").Value = "'6.511"
$ws.Range("E37").Value = "  -6.76%  "

$ws.Range("D38").Value = "'0.2628"
$ws.Range("E38").Value = "  -6.26%  "

$ws.Range("D39").Value = "'0.09122"
$ws.Range("E39").Value = "  -3.49%  "

$ws.Range("E40").Value = "  +0.66%  "

$ws.Range("D41").Value = "'13.30"
$ws.Range("E41").Value = "  -2.09%  "

$ws.Range("E42").Value = "  -6.54%  "

$ws.Range("D43").Value = "'0.7393"
$ws.Range("E43").Value = "  -6.29%  "

$ws.Range("D44").Value = "'15.72"
$ws.Range("E44").Value = "  -4.57%  "

$ws.Range("D45").Value = "'0.6782"
$ws.Range("E45").Value = "  -4.62%  "

$ws.Range("D46").Value = "'2.422"
$ws.Range("E46").Value = "  -4.94%  "

$ws.Range("D47").Value = "'4.047"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "'0.08204"
$ws.Range("E49").Value = "  -4.67%  "

$ws.Range("D50").Value = "'131.92"
$ws.Range("E50").Value = "  -4.00%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'83.05"
$ws.Range("E51").Value = "  -5.46%  "
